$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop LMU_NEURO_COGN_PSYCHOLOGY and FU_BERLIN_COGN_NEUROSCIENCE from the program
# list, shifting the remaining programs up two rows...
$ws.Range("A3").Value = "TU_BERLIN_COMP_NEUROSCIENCE"
$ws.Range("A4").Value = "UNI_BREMEN_NEUROSCIENCES"
$ws.Range("A5").Value = "UNI_OLDENBURG_NEUROSCIENCES"
$ws.Range("A6").Value = "UNI_OLDENBURG_NEUROCOGN_PSY"
# ...and add the new psy course analysis program where the list now ends
$ws.Range("A7").Value = "TU_DARMSTADT_COGNITIVE_SCIENCE"
$ws.Range("B7").Value = "Yes"

# The table is now one program shorter overall (7 -> 6), so the old last row is blank
$ws.Range("A8:B8").Clear()

# Shrink the "Yes/No" dropdown validation on column B to match the new data extent
$ws.Range("B8").Validation.Delete()

# The sheet grid shrinks by the two deleted rows as well
$ws.Rows("999:1000").Delete()
